$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: restyle A5/B5 and C5/D5/E5 to match the "continuation" look
# (same style family used by row 3) while keeping the existing values.
$ws.Range("A3:E3").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("B5").Value = 571

# --- Lay out the three new rows' formatting first (styles + row heights),
# copied from existing rows so the same style indices get reused.
$ws.Range("A2:E2").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Rows.Item(6).RowHeight = 43.2

$ws.Range("A2:E2").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A7").Clear()
$ws.Rows.Item(7).RowHeight = 21.6

$ws.Range("A2:E2").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A8").Clear()

# --- Now fill in the content column-by-column (this is the order the
# strings were actually typed in, and it determines shared-string order).
$ws.Range("C6").Value = ' Aaaah!'
$ws.Range("C7").Value = ' Eureka![K]\nIt sounds good to hear a shout like that, huh?'
$ws.Range("C8").Value = ' Please continue to explore!'

$ws.Range("A6").Value = "SCRIPT/P01P04A/us3111.ssb"

$ws.Range("D6").Value = ' Аааах!'
$ws.Range("D7").Value = ' Эврика![K] Приятно ведь слышать\nтакое слово, да?'
$ws.Range("D8").Value = ' Прошу, продолжайте исследовать!'

$ws.Range("E6").Value = ' Ààààö!'
$ws.Range("E7").Value = ' Üâñéëà![K] Ðñéÿóîï âåäû òìúšàóû\nóàëïå òìïâï, äà?'
$ws.Range("E8").Value = ' Ðñïšô, ðñïäïìçàêóå éòòìåäïâàóû!'

$ws.Range("B6").Value = 533
$ws.Range("B7").Value = 536
$ws.Range("B8").Value = 547

# --- Update the view: select A6, matching the saved UI state after the new
# rows were added and scrolled into view.
$ws.Range("A6").Select()
